$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '76.123.84'
$ws.Range("E2").Value = '  +1.00%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.917.17'
$ws.Range("E3").Value = '  +5.39%  '

# Row 4
$ws.Range("E4").Value = '  +0.04%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '199.24'
$ws.Range("E5").Value = '  +6.01%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '605.81'
$ws.Range("E6").Value = '  +1.42%  '

# Row 7
$ws.Range("E7").Value = '  +0.01%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.556'
$ws.Range("E8").Value = '  +3.04%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.195'
$ws.Range("E9").Value = '  -2.96%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '2.917.02'
$ws.Range("E10").Value = '  +5.38%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.416'
$ws.Range("E11").Value = '  +13.54%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.161'
$ws.Range("E12").Value = '  -0.71%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.96'
$ws.Range("E13").Value = '  +3.46%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.488.90'
$ws.Range("E14").Value = '  +6.95%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '75.973.58'
$ws.Range("E15").Value = '  +0.94%  '

# Row 16
$ws.Range("B16").Value = 'Avalanche'
$ws.Range("C16").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '27.93'
$ws.Range("E16").Value = '  +1.81%  '

# Row 17
$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000192'
$ws.Range("E17").Value = '  +0.11%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.928.06'
$ws.Range("E18").Value = '  +6.56%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '8.93'
$ws.Range("E19").Value = '  -4.27%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.76'
$ws.Range("E20").Value = '  +4.06%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '382.91'
$ws.Range("E21").Value = '  +0.50%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.35'
$ws.Range("E22").Value = '  +1.65%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.20'
$ws.Range("E23").Value = '  +1.90%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '72.40'
$ws.Range("E24").Value = '  +1.74%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.999'
$ws.Range("E25").Value = '  -0.04%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.062.50'
$ws.Range("E26").Value = '  +6.57%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '4.31'
$ws.Range("E27").Value = '  +1.21%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.86'
$ws.Range("E28").Value = '  +2.83%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0000112'
$ws.Range("E29").Value = '  +9.14%  '

# Row 30
$ws.Range("E30").Value = '  +0.14%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.43'
$ws.Range("E31").Value = '  +0.79%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '514.29'
$ws.Range("E32").Value = '  -3.59%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.88'
$ws.Range("E33").Value = '  -0.91%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.85'
$ws.Range("E34").Value = '  +2.08%  '

# Row 35
$ws.Range("E35").Value = '  +0.06%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '164.42'
$ws.Range("E36").Value = '  +1.43%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '20.23'
$ws.Range("E37").Value = '  +1.12%  '

# Row 38
$ws.Range("B38").Value = 'Kaspa'
$ws.Range("C38").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.116'
$ws.Range("E38").Value = '  -3.59%  '

# Row 39
$ws.Range("B39").Value = 'WhiteBITCoin'
$ws.Range("C39").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '19.71'
$ws.Range("E39").Value = '  +1.92%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '182.49'
$ws.Range("E40").Value = '  +1.38%  '

# Row 41
$ws.Range("E41").Value = '  +0.04%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.348'
$ws.Range("E42").Value = '  +3.38%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.11'
$ws.Range("E43").Value = '  +0.13%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.70'
$ws.Range("E44").Value = '  -0.73%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0930'
$ws.Range("E45").Value = '  +7.42%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.24'
$ws.Range("E46").Value = '  +0.23%  '

# Row 47
$ws.Range("B47").Value = 'dogwifhat'
$ws.Range("C47").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.47'
$ws.Range("E47").Value = '  +1.57%  '

# Row 48
$ws.Range("B48").Value = 'OKB'
$ws.Range("C48").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '40.47'
$ws.Range("E48").Value = '  +2.32%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.585'
$ws.Range("E49").Value = '  +5.18%  '

# Row 50
$ws.Range("B50").Value = 'Filecoin'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.79'
$ws.Range("E50").Value = '  +1.81%  '

# Row 51
$ws.Range("B51").Value = 'Mantle'
$ws.Range("C51").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.661'
$ws.Range("E51").Value = '  +9.19%  '
